$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Rows 91 and 92 got their betting-odds payloads swapped (id/div/date
# columns stay put, everything from the match-result column onward
# trades places between the two rows).
# ---------------------------------------------------------------------
$row91 = [ordered]@{
    B = 6924568
    F = "Atletico Morelia"
    G = "Atlante"
    H = 0
    J = "A"
    K = 2.4
    L = 3
    M = 2.875
    N = 2.7
    O = 3.1
    P = 2.8
    Q = 0
    R = 1.85
    S = 1.95
    T = 2.25
    U = 1.975
    V = 1.725
    W = -1
    X = -1
    Y = 1.8
    Z = -1
    AA = 0.95
    AB = -1
    AC = 0.7250000000000001
}

$row92 = [ordered]@{
    B = 6924569
    F = "Venados FC"
    G = "Dorados"
    H = 4
    J = "H"
    K = 1.615
    L = 4
    M = 4.5
    N = 1.5
    O = 4.75
    P = 5.75
    Q = -1.25
    R = 1.925
    S = 1.875
    T = 3
    U = 1.75
    V = 1.95
    W = 0.5
    X = -1
    Y = -1
    Z = 0.925
    AA = -1
    AB = 0.75
    AC = -1
}

foreach ($col in $row92.Keys) {
    $ws.Range("$col" + "91").Value = $row92[$col]
}
foreach ($col in $row91.Keys) {
    $ws.Range("$col" + "92").Value = $row91[$col]
}

# ---------------------------------------------------------------------
# New row 225: an upcoming (not-yet-played) Universidad Guadalajara vs
# Tepatitlan FC fixture appended to the bottom of the data set.
# ---------------------------------------------------------------------
$ws.Range("A224").Copy($ws.Range("A225")) | Out-Null
$ws.Range("E224").Copy($ws.Range("E225")) | Out-Null

$ws.Range("A225").Value = 223
$ws.Range("B225").Value = 7641726
$ws.Range("C225").Value = "Mexico Liga de Expansion"
$ws.Range("D225").Value = "Mexico Liga de Expansion"
$ws.Range("E225").Value = 45393.00347222222
$ws.Range("F225").Value = "Universidad Guadalajara"
$ws.Range("G225").Value = "Tepatitlan FC"
$ws.Range("K225").Value = 1.25
$ws.Range("L225").Value = 5.5
$ws.Range("M225").Value = 7.5
$ws.Range("N225").Value = 1.3
$ws.Range("O225").Value = 5.25
$ws.Range("P225").Value = 6.5
$ws.Range("Q225").Value = -1.5
$ws.Range("R225").Value = 1.9
$ws.Range("S225").Value = 1.9
$ws.Range("T225").Value = 2.75
$ws.Range("U225").Value = 1.875
$ws.Range("V225").Value = 1.925
$ws.Range("W225").Value = 0
$ws.Range("X225").Value = 0
$ws.Range("Y225").Value = 0
$ws.Range("Z225").Value = 0
$ws.Range("AA225").Value = 0
